$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff: new crypto price/volume snapshot + one coin swap (row 51)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.380.59'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.22%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.583.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.06'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.47%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.492'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.83'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.60%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.92%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0894'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.35%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.807.82'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.584.18'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.11%  '

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.43%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '28.415.89'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.24%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.01'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.91%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.82'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.49%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0688'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.61%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.52%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.05'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.08%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.34'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.35%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.48%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.76%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0484'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.84%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.20'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.74%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.54%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.401.28'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.86%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.09'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.60%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.36'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.31%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.47%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0164'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.79%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.520'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.62%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.13%  '

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.80%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0460'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.51%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.73%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.74'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.35%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.928'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.38%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.719.74'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.38'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.42%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0517'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.16%  '
